$d = $word.ActiveDocument

$d.Content.Find.Execute("728×9=6552", $true, $false, $false, $false, $false, $true, 1, $false, "677×9=6093", 2) | Out-Null
$d.Content.Find.Execute("806×7=5642", $true, $false, $false, $false, $false, $true, 1, $false, "893×6=5358", 2) | Out-Null
$d.Content.Find.Execute("295×9=2655", $true, $false, $false, $false, $false, $true, 1, $false, "168×5=840", 2) | Out-Null
$d.Content.Find.Execute("228×8=1824", $true, $false, $false, $false, $false, $true, 1, $false, "407×3=1221", 2) | Out-Null
$d.Content.Find.Execute("734×2=1468", $true, $false, $false, $false, $false, $true, 1, $false, "107×9=963", 2) | Out-Null
$d.Content.Find.Execute("407×6=2442", $true, $false, $false, $false, $false, $true, 1, $false, "712×6=4272", 2) | Out-Null
$d.Content.Find.Execute("149×6=894", $true, $false, $false, $false, $false, $true, 1, $false, "101×4=404", 2) | Out-Null
$d.Content.Find.Execute("895×8=7160", $true, $false, $false, $false, $false, $true, 1, $false, "411×3=1233", 2) | Out-Null
$d.Content.Find.Execute("469×3=1407", $true, $false, $false, $false, $false, $true, 1, $false, "587×8=4696", 2) | Out-Null
$d.Content.Find.Execute("812×6=4872", $true, $false, $false, $false, $false, $true, 1, $false, "722×6=4332", 2) | Out-Null
$d.Content.Find.Execute("467×5=2335", $true, $false, $false, $false, $false, $true, 1, $false, "644×8=5152", 2) | Out-Null
$d.Content.Find.Execute("675×8=5400", $true, $false, $false, $false, $false, $true, 1, $false, "430×6=2580", 2) | Out-Null
$d.Content.Find.Execute("712×7=4984", $true, $false, $false, $false, $false, $true, 1, $false, "292×4=1168", 2) | Out-Null
$d.Content.Find.Execute("621×6=3726", $true, $false, $false, $false, $false, $true, 1, $false, "823×7=5761", 2) | Out-Null
$d.Content.Find.Execute("840×3=2520", $true, $false, $false, $false, $false, $true, 1, $false, "660×6=3960", 2) | Out-Null
$d.Content.Find.Execute("579×4=2316", $true, $false, $false, $false, $false, $true, 1, $false, "367×3=1101", 2) | Out-Null
$d.Content.Find.Execute("688×5=3440", $true, $false, $false, $false, $false, $true, 1, $false, "262×8=2096", 2) | Out-Null
$d.Content.Find.Execute("968×7=6776", $true, $false, $false, $false, $false, $true, 1, $false, "735×3=2205", 2) | Out-Null
$d.Content.Find.Execute("119×2=238", $true, $false, $false, $false, $false, $true, 1, $false, "342×8=2736", 2) | Out-Null
$d.Content.Find.Execute("238×7=1666", $true, $false, $false, $false, $false, $true, 1, $false, "244×2=488", 2) | Out-Null
$d.Content.Find.Execute("120×6=720", $true, $false, $false, $false, $false, $true, 1, $false, "854×9=7686", 2) | Out-Null
$d.Content.Find.Execute("944×2=1888", $true, $false, $false, $false, $false, $true, 1, $false, "770×8=6160", 2) | Out-Null
$d.Content.Find.Execute("651×5=3255", $true, $false, $false, $false, $false, $true, 1, $false, "401×7=2807", 2) | Out-Null
$d.Content.Find.Execute("518×5=2590", $true, $false, $false, $false, $false, $true, 1, $false, "939×3=2817", 2) | Out-Null
$d.Content.Find.Execute("172×2=344", $true, $false, $false, $false, $false, $true, 1, $false, "664×9=5976", 2) | Out-Null
